$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" "63.681.23"
Set-TextCell "E2" "  +2.79%  "
Set-TextCell "D3" "3.132.02"
Set-TextCell "E3" "  +1.58%  "
Set-TextCell "E4" "  +0.05%  "
Set-TextCell "D5" "588.47"
Set-TextCell "E5" "  +1.47%  "
Set-TextCell "D6" "146.49"
Set-TextCell "E6" "  +2.83%  "
Set-TextCell "E7" "  +0.04%  "
Set-TextCell "D8" "3.123.07"
Set-TextCell "E8" "  +1.72%  "
Set-TextCell "D9" "0.533"
Set-TextCell "E9" "  +1.22%  "
Set-TextCell "D10" "0.160"
Set-TextCell "E10" "  +13.98%  "
Set-TextCell "D11" "5.71"
Set-TextCell "E11" "  +2.25%  "
Set-TextCell "D12" "0.470"
Set-TextCell "E12" "  +0.54%  "
Set-TextCell "D13" "0.0000251"
Set-TextCell "E13" "  +4.65%  "
Set-TextCell "D14" "36.53"
Set-TextCell "E14" "  +3.35%  "
Set-TextCell "E15" "  -0.69%  "
Set-TextCell "D16" "3.659.40"
Set-TextCell "E16" "  +1.97%  "
Set-TextCell "D17" "7.17"
Set-TextCell "E17" "  -1.48%  "
Set-TextCell "D18" "63.638.15"
Set-TextCell "E18" "  +2.90%  "
Set-TextCell "D19" "3.133.70"
Set-TextCell "E19" "  +1.97%  "
Set-TextCell "D20" "463.82"
Set-TextCell "E20" "  +3.16%  "
Set-TextCell "D21" "14.43"
Set-TextCell "E21" "  +3.71%  "
Set-TextCell "D22" "0.734"
Set-TextCell "E22" "  +0.68%  "
Set-TextCell "D23" "7.54"
Set-TextCell "E23" "  +1.45%  "
Set-TextCell "D24" "13.24"
Set-TextCell "E24" "  -4.00%  "
Set-TextCell "D25" "82.24"
Set-TextCell "E25" "  +0.35%  "
Set-TextCell "E26" "  -0.06%  "
Set-TextCell "D27" "8.94"
Set-TextCell "E27" "  +9.89%  "
Set-TextCell "D28" "2.70"
Set-TextCell "E28" "  +1.60%  "
Set-TextCell "E29" "  -1.79%  "
Set-TextCell "E30" "  +0.05%  "
Set-TextCell "D31" "6.89"
Set-TextCell "E31" "  +1.85%  "
Set-TextCell "D32" "27.12"
Set-TextCell "E32" "  +1.38%  "
Set-TextCell "D33" "0.109"
Set-TextCell "E33" "  -1.79%  "
Set-TextCell "D34" "0.0₃0866"
Set-TextCell "E34" "  +7.91%  "
Set-TextCell "D35" "2.37"
Set-TextCell "E35" "  +7.71%  "
Set-TextCell "D36" "1.05"
Set-TextCell "E36" "  +1.32%  "
Set-TextCell "D37" "3.36"
Set-TextCell "E37" "  +11.75%  "
Set-TextCell "D38" "6.09"
Set-TextCell "E38" "  +0.57%  "
Set-TextCell "D39" "51.01"
Set-TextCell "E39" "  +1.15%  "
Set-TextCell "D40" "447.56"
Set-TextCell "E40" "  +3.97%  "
Set-TextCell "D41" "8.75"
Set-TextCell "E41" "  -0.98%  "
Set-TextCell "E42" "  +0.05%  "
Set-TextCell "D43" "2.893.82"
Set-TextCell "E43" "  +3.39%  "
Set-TextCell "E44" "  +2.77%  "
Set-TextCell "E45" "  +2.50%  "
Set-TextCell "D46" "2.18"
Set-TextCell "E46" "  +3.22%  "
Set-TextCell "D47" "36.39"
Set-TextCell "E47" "  +3.23%  "
Set-TextCell "D48" "124.85"
Set-TextCell "E48" "  +0.90%  "
Set-TextCell "E49" "  +0.06%  "
Set-TextCell "D50" "0.111"
Set-TextCell "E50" "  +0.13%  "
Set-TextCell "D51" "24.68"
Set-TextCell "E51" "  +2.40%  "
